$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "99-2=97"
$t.Cell(1,2).Range.Text = "87-44=43"
$t.Cell(1,3).Range.Text = "28+43=71"
$t.Cell(1,4).Range.Text = "71+8=79"
$t.Cell(1,5).Range.Text = "47+10=57"
$t.Cell(2,1).Range.Text = "97-60=37"
$t.Cell(2,2).Range.Text = "6+75=81"
$t.Cell(2,3).Range.Text = "4+6=10"
$t.Cell(2,4).Range.Text = "65-31=34"
$t.Cell(2,5).Range.Text = "75+11=86"
$t.Cell(3,1).Range.Text = "96-90=6"
$t.Cell(3,2).Range.Text = "25+61=86"
$t.Cell(3,3).Range.Text = "64-22=42"
$t.Cell(3,4).Range.Text = "95-4=91"
$t.Cell(3,5).Range.Text = "7+60=67"
$t.Cell(4,1).Range.Text = "53-37=16"
$t.Cell(4,2).Range.Text = "38+52=90"
$t.Cell(4,3).Range.Text = "50+19=69"
$t.Cell(4,4).Range.Text = "51+46=97"
$t.Cell(4,5).Range.Text = "85-42=43"
$t.Cell(5,1).Range.Text = "39+3=42"
$t.Cell(5,2).Range.Text = "68+3=71"
$t.Cell(5,3).Range.Text = "80+6=86"
$t.Cell(5,4).Range.Text = "71-41=30"
$t.Cell(5,5).Range.Text = "2+16=18"
$t.Cell(6,1).Range.Text = "33-33=0"
$t.Cell(6,2).Range.Text = "12+68=80"
$t.Cell(6,3).Range.Text = "34-23=11"
$t.Cell(6,4).Range.Text = "26+63=89"
$t.Cell(6,5).Range.Text = "35+24=59"
$t.Cell(7,1).Range.Text = "33-25=8"
$t.Cell(7,2).Range.Text = "45-8=37"
$t.Cell(7,3).Range.Text = "70+20=90"
$t.Cell(7,4).Range.Text = "68+0=68"
$t.Cell(7,5).Range.Text = "20-13=7"
$t.Cell(8,1).Range.Text = "39+58=97"
$t.Cell(8,2).Range.Text = "97-20=77"
$t.Cell(8,3).Range.Text = "30-12=18"
$t.Cell(8,4).Range.Text = "75+8=83"
$t.Cell(8,5).Range.Text = "9+66=75"
$t.Cell(9,1).Range.Text = "16+39=55"
$t.Cell(9,2).Range.Text = "3+41=44"
$t.Cell(9,3).Range.Text = "4+95=99"
$t.Cell(9,4).Range.Text = "73+18=91"
$t.Cell(9,5).Range.Text = "63+26=89"
$t.Cell(10,1).Range.Text = "8+67=75"
$t.Cell(10,2).Range.Text = "49+25=74"
$t.Cell(10,3).Range.Text = "28+52=80"
$t.Cell(10,4).Range.Text = "54+7=61"
$t.Cell(10,5).Range.Text = "52+45=97"
$t.Cell(11,1).Range.Text = "30-10=20"
$t.Cell(11,2).Range.Text = "71-34=37"
$t.Cell(11,3).Range.Text = "80-22=58"
$t.Cell(11,4).Range.Text = "22+55=77"
$t.Cell(11,5).Range.Text = "54-40=14"
$t.Cell(12,1).Range.Text = "15-6=9"
$t.Cell(12,2).Range.Text = "14+41=55"
$t.Cell(12,3).Range.Text = "88-54=34"
$t.Cell(12,4).Range.Text = "93-4=89"
$t.Cell(12,5).Range.Text = "50+0=50"
$t.Cell(13,1).Range.Text = "92-47=45"
$t.Cell(13,2).Range.Text = "6+48=54"
$t.Cell(13,3).Range.Text = "19-9=10"
$t.Cell(13,4).Range.Text = "90-4=86"
$t.Cell(13,5).Range.Text = "82+3=85"
$t.Cell(14,1).Range.Text = "15+45=60"
$t.Cell(14,2).Range.Text = "58-14=44"
$t.Cell(14,3).Range.Text = "30-24=6"
$t.Cell(14,4).Range.Text = "99-76=23"
$t.Cell(14,5).Range.Text = "10+5=15"
$t.Cell(15,1).Range.Text = "26-19=7"
$t.Cell(15,2).Range.Text = "18+2=20"
$t.Cell(15,3).Range.Text = "89+0=89"
$t.Cell(15,4).Range.Text = "99-64=35"
$t.Cell(15,5).Range.Text = "7+8=15"
$t.Cell(16,1).Range.Text = "31-24=7"
$t.Cell(16,2).Range.Text = "78-27=51"
$t.Cell(16,3).Range.Text = "12+79=91"
$t.Cell(16,4).Range.Text = "71-5=66"
$t.Cell(16,5).Range.Text = "38+15=53"
$t.Cell(17,1).Range.Text = "43+45=88"
$t.Cell(17,2).Range.Text = "92-40=52"
$t.Cell(17,3).Range.Text = "7+51=58"
$t.Cell(17,4).Range.Text = "35-5=30"
$t.Cell(17,5).Range.Text = "14+11=25"
$t.Cell(18,1).Range.Text = "37+47=84"
$t.Cell(18,2).Range.Text = "25+64=89"
$t.Cell(18,3).Range.Text = "13+61=74"
$t.Cell(18,4).Range.Text = "16+24=40"
$t.Cell(18,5).Range.Text = "68-1=67"
$t.Cell(19,1).Range.Text = "51+24=75"
$t.Cell(19,2).Range.Text = "24+10=34"
$t.Cell(19,3).Range.Text = "46-25=21"
$t.Cell(19,4).Range.Text = "97-75=22"
$t.Cell(19,5).Range.Text = "57-26=31"
$t.Cell(20,1).Range.Text = "28+14=42"
$t.Cell(20,2).Range.Text = "38+51=89"
$t.Cell(20,3).Range.Text = "6+64=70"
$t.Cell(20,4).Range.Text = "24+48=72"
$t.Cell(20,5).Range.Text = "79-20=59"

Write-Output "done"
